# Auto-generated files on 2025-10-09
# Update HotStock_Top20 sheet values (columns A-C, rows 2-21) to reflect
# the refreshed "hot stocks" ranking lists.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @("通富微电", "上海电气", "山子高科")
    3  = @("北方稀土", "山子高科", "海南华铁")
    4  = @("深科技",   "北方稀土", "北方稀土")
    5  = @("山子高科", "通富微电", "赣锋锂业")
    6  = @("中芯国际", "赣锋锂业", "上海电气")
    7  = @("上海电气", "紫金矿业", "中芯国际")
    8  = @("合锻智能", "江西铜业", "紫金矿业")
    9  = @("赣锋锂业", "特变电工", "深科技")
    10 = @("紫金矿业", "深科技",   "通富微电")
    11 = @("领益智造", "永鼎股份", "赛力斯")
    12 = @("永鼎股份", "中兴通讯", "蓝丰生化")
    13 = @("洛阳钼业", "合锻智能", "天际股份")
    14 = @("江西铜业", "洛阳钼业", "云汉芯城")
    15 = @("海南华铁", "中芯国际", "张江高科")
    16 = @("中兴通讯", "东方财富", "洛阳钼业")
    17 = @("特变电工", "北方铜业", "三花智控")
    18 = @("张江高科", "领益智造", "金力永磁")
    19 = @("三花智控", "三花智控", "万向钱潮")
    20 = @("天际股份", "张江高科", "领益智造")
    21 = @("海鸥住工", "融发核电", "合锻智能")
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
}
